# Auto commit at 2025-11-09 10:04:41.97
# Updates the "Metrics" sheet values (B2:B13) with refreshed figures; the
# "today" sheet pulls these via formulas (=Metrics!Bn) and will recalc
# automatically. Also restores the saved selection on each sheet.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 104979.72
$metrics.Range("B3").Value = 91628.68
$metrics.Range("B4").Value = 32663.47
$metrics.Range("B5").Value = 4413
$metrics.Range("B6").Value = 4901225.47
$metrics.Range("B7").Value = 4133705.3600000008
$metrics.Range("B8").Value = 1439623.2999999998
$metrics.Range("B9").Value = 190620
$metrics.Range("B10").Value = 33366606.460000005
$metrics.Range("B11").Value = 31408980.52
$metrics.Range("B12").Value = 11721345.340000004
$metrics.Range("B13").Value = 1288250

$metrics.Activate()
$metrics.Range("D13").Select()

$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("C5").Select()

$excel.Calculate()
